{"js": "// Replace every math-problem cell in the (single) table with the new\n// set of values, preserving row/column layout, cell formatting and\n// paragraph/run properties (Word keeps those when Range.Text /\n// table.values rewrites only the text).\nconst table = context.document.body.tables.getFirst();\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst newValues = [\n  [\"4+37=\", \"9+16=\", \"82-55=\", \"27+36=\", \"43+29=\"],\n  [\"80-47=\", \"30-8=\", \"25+4=\", \"42-1=\", \"53-18=\"],\n  [\"28+3=\", \"87-34=\", \"22+18=\", \"25+17=\", \"18+57=\"],\n  [\"36-16=\", \"28+59=\", \"18-13=\", \"72-63=\", \"38+30=\"],\n  [\"3+90=\", \"35-25=\", \"35-6=\", \"25+36=\", \"40+29=\"],\n  [\"8+32=\", \"72+0=\", \"83-11=\", \"44-25=\", \"85-17=\"],\n  [\"64+22=\", \"83-76=\", \"6+16=\", \"70-30=\", \"35+6=\"],\n  [\"24+44=\", \"59+12=\", \"43-25=\", \"96-39=\", \"79-25=\"],\n  [\"54+13=\", \"23+29=\", \"98-48=\", \"32+46=\", \"40+55=\"],\n  [\"28+21=\", \"83-34=\", \"39+43=\", \"47+6=\", \"86-17=\"],\n  [\"30+53=\", \"61+36=\", \"72-33=\", \"73+25=\", \"37+14=\"],\n  [\"74-1=\", \"20+27=\", \"55+26=\", \"81-52=\", \"42+37=\"],\n  [\"37+54=\", \"81-27=\", \"73-47=\", \"58+6=\", \"49-13=\"],\n  [\"77-38=\", \"89-7=\", \"75-17=\", \"36+3=\", \"53-42=\"],\n  [\"42+1=\", \"23+1=\", \"47+24=\", \"36+54=\", \"72+19=\"],\n  [\"68+26=\", \"92-57=\", \"20+39=\", \"33+66=\", \"96-63=\"],\n  [\"6+68=\", \"91-22=\", \"53+28=\", \"95-34=\", \"23+76=\"],\n  [\"19-16=\", \"75-9=\", \"98-20=\", \"3+24=\", \"30+15=\"],\n  [\"2+29=\", \"43-30=\", \"15+66=\", \"59-31=\", \"74-2=\"],\n  [\"66+9=\", \"0+87=\", \"25+50=\", \"92-54=\", \"99-63=\"]\n];\n\nconst currentColumnCount = table.values[0].length;\nif (table.rowCount !== newValues.length || currentColumnCount !== newValues[0].length) {\n  throw new Error(\n    `Unexpected table shape: ${table.rowCount}x${currentColumnCount}, ` +\n    `expected ${newValues.length}x${newValues[0].length}`\n  );\n}\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Replace every math-problem cell in the (single) table with the new\n# set of values, preserving row/column layout and cell formatting\n# (Range.Text only rewrites the text run inside the existing cell mark).\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$newValues = @(\n    @(\"4+37=\", \"9+16=\", \"82-55=\", \"27+36=\", \"43+29=\"),\n    @(\"80-47=\", \"30-8=\", \"25+4=\", \"42-1=\", \"53-18=\"),\n    @(\"28+3=\", \"87-34=\", \"22+18=\", \"25+17=\", \"18+57=\"),\n    @(\"36-16=\", \"28+59=\", \"18-13=\", \"72-63=\", \"38+30=\"),\n    @(\"3+90=\", \"35-25=\", \"35-6=\", \"25+36=\", \"40+29=\"),\n    @(\"8+32=\", \"72+0=\", \"83-11=\", \"44-25=\", \"85-17=\"),\n    @(\"64+22=\", \"83-76=\", \"6+16=\", \"70-30=\", \"35+6=\"),\n    @(\"24+44=\", \"59+12=\", \"43-25=\", \"96-39=\", \"79-25=\"),\n    @(\"54+13=\", \"23+29=\", \"98-48=\", \"32+46=\", \"40+55=\"),\n    @(\"28+21=\", \"83-34=\", \"39+43=\", \"47+6=\", \"86-17=\"),\n    @(\"30+53=\", \"61+36=\", \"72-33=\", \"73+25=\", \"37+14=\"),\n    @(\"74-1=\", \"20+27=\", \"55+26=\", \"81-52=\", \"42+37=\"),\n    @(\"37+54=\", \"81-27=\", \"73-47=\", \"58+6=\", \"49-13=\"),\n    @(\"77-38=\", \"89-7=\", \"75-17=\", \"36+3=\", \"53-42=\"),\n    @(\"42+1=\", \"23+1=\", \"47+24=\", \"36+54=\", \"72+19=\"),\n    @(\"68+26=\", \"92-57=\", \"20+39=\", \"33+66=\", \"96-63=\"),\n    @(\"6+68=\", \"91-22=\", \"53+28=\", \"95-34=\", \"23+76=\"),\n    @(\"19-16=\", \"75-9=\", \"98-20=\", \"3+24=\", \"30+15=\"),\n    @(\"2+29=\", \"43-30=\", \"15+66=\", \"59-31=\", \"74-2=\"),\n    @(\"66+9=\", \"0+87=\", \"25+50=\", \"92-54=\", \"99-63=\")\n)\n\nif ($tbl.Rows.Count -ne $newValues.Count -or $tbl.Columns.Count -ne $newValues[0].Count) {\n    throw \"Unexpected table shape: $($tbl.Rows.Count)x$($tbl.Columns.Count), expected $($newValues.Count)x$($newValues[0].Count)\"\n}\n\nfor ($r = 1; $r -le $tbl.Rows.Count; $r++) {\n    for ($c = 1; $c -le $tbl.Columns.Count; $c++) {\n        $tbl.Cell($r, $c).Range.Text = $newValues[$r - 1][$c - 1]\n    }\n}\n"}
